$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Date heading
Replace-Text "2023-08-24 Thursday" "2023-08-25 Friday"

# Table row 1 (unique texts -> safe to use Find/Replace)
Replace-Text "56÷5=11, 1" "47÷8=5, 7"
Replace-Text "20÷8=2, 4" "36÷5=7, 1"
Replace-Text "25÷8=3, 1" "50÷2=25, 0"
Replace-Text "79÷4=19, 3" "83÷9=9, 2"

# Table row 5
Replace-Text "75÷5=15, 0" "66÷2=33, 0"
Replace-Text "32÷7=4, 4" "95÷7=13, 4"
Replace-Text "87÷4=21, 3" "49÷4=12, 1"
Replace-Text "67÷9=7, 4" "40÷3=13, 1"

# Table row 9
Replace-Text "34÷8=4, 2" "48÷4=12, 0"
Replace-Text "67÷3=22, 1" "14÷4=3, 2"
Replace-Text "92÷3=30, 2" "69÷7=9, 6"
Replace-Text "87÷3=29, 0" "81÷8=10, 1"
Replace-Text "98÷9=10, 8" "39÷2=19, 1"

# Table row 13
Replace-Text "40÷8=5, 0" "24÷9=2, 6"
Replace-Text "60÷5=12, 0" "77÷6=12, 5"
Replace-Text "91÷7=13, 0" "93÷3=31, 0"
Replace-Text "84÷5=16, 4" "46÷5=9, 1"
Replace-Text "26÷5=5, 1" "85÷7=12, 1"

# Table row 17
Replace-Text "39÷8=4, 7" "69÷2=34, 1"
Replace-Text "89÷5=17, 4" "28÷9=3, 1"
Replace-Text "80÷9=8, 8" "19÷5=3, 4"
Replace-Text "28÷6=4, 4" "25÷3=8, 1"
Replace-Text "91÷5=18, 1" "15÷4=3, 3"

# "79÷5=15, 4" appears twice (row 1 col 5, row 5 col 5) with different
# replacements, so Find/Replace can't disambiguate them - address the
# cells directly via the table object model instead.
$t = $d.Tables.Item(1)
$t.Cell(1, 5).Range.Text = "38÷9=4, 2"
$t.Cell(5, 5).Range.Text = "94÷8=11, 6"
